$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 204
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A204:L204").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("N204").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(204,1).Value = 45190.40807770833
$ws.Cells.Item(204,2).Value = "drmiso526@naver.com"
$ws.Cells.Item(204,3).Value = "사회복지학과"
$ws.Cells.Item(204,4).Value = 20192310
$ws.Cells.Item(204,5).Value = "김세중"
$ws.Cells.Item(204,6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(204,7).Value = 0.1
$ws.Cells.Item(204,8).Value = "6:4"
$ws.Cells.Item(204,9).Value = "15분의 1"
$ws.Cells.Item(204,10).Value = "44만호, 153만명"
$ws.Cells.Item(204,11).Value = "충청"
$ws.Cells.Item(204,12).Value = "Black"
$ws.Cells.Item(204,14).Value = "찬성한다."

# Row 205
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A205:L205").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("N205").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(205,1).Value = 45190.41620032408
$ws.Cells.Item(205,2).Value = "chlwnsgur1202@gmail.com"
$ws.Cells.Item(205,3).Value = "인공지능융합학부"
$ws.Cells.Item(205,4).Value = 20236783
$ws.Cells.Item(205,5).Value = "최준혁"
$ws.Cells.Item(205,6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(205,7).Value = 0.1
$ws.Cells.Item(205,8).Value = "6:4"
$ws.Cells.Item(205,9).Value = "20분의 1"
$ws.Cells.Item(205,10).Value = "20만호, 69만명"
$ws.Cells.Item(205,11).Value = "경상"
$ws.Cells.Item(205,12).Value = "Black"
$ws.Cells.Item(205,14).Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# Row 206
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A206:L206").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("M206").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(206,1).Value = 45190.45822483796
$ws.Cells.Item(206,2).Value = "rjsgjsd135@naver.com"
$ws.Cells.Item(206,3).Value = "사회복지학과"
$ws.Cells.Item(206,4).Value = 20181077
$ws.Cells.Item(206,5).Value = "이은혁"
$ws.Cells.Item(206,6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(206,7).Value = 0.1
$ws.Cells.Item(206,8).Value = "6:4"
$ws.Cells.Item(206,9).Value = "30분의 1"
$ws.Cells.Item(206,10).Value = "20만호, 69만명"
$ws.Cells.Item(206,11).Value = "충청"
$ws.Cells.Item(206,12).Value = "Red"
$ws.Cells.Item(206,13).Value = "모름/무응답"

# Row 207
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A207:L207").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("N207").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(207,1).Value = 45190.47622938658
$ws.Cells.Item(207,2).Value = "lsk8424@naver.com"
$ws.Cells.Item(207,3).Value = "정치행정학과"
$ws.Cells.Item(207,4).Value = 20172428
$ws.Cells.Item(207,5).Value = "이상기"
$ws.Cells.Item(207,6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(207,7).Value = 0.1
$ws.Cells.Item(207,8).Value = "6:4"
$ws.Cells.Item(207,9).Value = "10분의 1"
$ws.Cells.Item(207,10).Value = "20만호, 69만명"
$ws.Cells.Item(207,11).Value = "충청"
$ws.Cells.Item(207,12).Value = "Black"
$ws.Cells.Item(207,14).Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# Row 208
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A208:L208").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("M208").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(208,1).Value = 45190.49009452546
$ws.Cells.Item(208,2).Value = "bravemw2@naver.com"
$ws.Cells.Item(208,3).Value = "데이터테크"
$ws.Cells.Item(208,4).Value = 20203257
$ws.Cells.Item(208,5).Value = "태민우"
$ws.Cells.Item(208,6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(208,7).Value = 0.1
$ws.Cells.Item(208,8).Value = "6:4"
$ws.Cells.Item(208,9).Value = "10분의 1"
$ws.Cells.Item(208,10).Value = "20만호, 69만명"
$ws.Cells.Item(208,11).Value = "평안"
$ws.Cells.Item(208,12).Value = "Red"
$ws.Cells.Item(208,13).Value = "반대한다."

# Row 209
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A209:L209").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("N209").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(209,1).Value = 45190.5015333912
$ws.Cells.Item(209,2).Value = "seraphic_0913@naver.com"
$ws.Cells.Item(209,3).Value = "법학과"
$ws.Cells.Item(209,4).Value = 20232720
$ws.Cells.Item(209,5).Value = "배윤서"
$ws.Cells.Item(209,6).Value = "실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다."
$ws.Cells.Item(209,7).Value = 0.3
$ws.Cells.Item(209,8).Value = "4:6"
$ws.Cells.Item(209,9).Value = "15분의 1"
$ws.Cells.Item(209,10).Value = "130만호, 5백만명"
$ws.Cells.Item(209,11).Value = "경기"
$ws.Cells.Item(209,12).Value = "Black"
$ws.Cells.Item(209,14).Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# Row 210
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A210:L210").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("M210").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(210,1).Value = 45190.51627489583
$ws.Cells.Item(210,2).Value = "cmsgood99@gmail.com"
$ws.Cells.Item(210,3).Value = "언어청각학부"
$ws.Cells.Item(210,4).Value = 20233905
$ws.Cells.Item(210,5).Value = "곽동희"
$ws.Cells.Item(210,6).Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Cells.Item(210,7).Value = 0.1
$ws.Cells.Item(210,8).Value = "3:7"
$ws.Cells.Item(210,9).Value = "20분의 1"
$ws.Cells.Item(210,10).Value = "130만호, 5백만명"
$ws.Cells.Item(210,11).Value = "충청"
$ws.Cells.Item(210,12).Value = "Red"
$ws.Cells.Item(210,13).Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# Row 211
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A211:L211").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("M211").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(211,1).Value = 45190.51763962963
$ws.Cells.Item(211,2).Value = "679qmg@gmail.com"
$ws.Cells.Item(211,3).Value = "법학과"
$ws.Cells.Item(211,4).Value = 20222750
$ws.Cells.Item(211,5).Value = "이혜원"
$ws.Cells.Item(211,6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(211,7).Value = 0.1
$ws.Cells.Item(211,8).Value = "6:4"
$ws.Cells.Item(211,9).Value = "20분의 1"
$ws.Cells.Item(211,10).Value = "20만호, 69만명"
$ws.Cells.Item(211,11).Value = "충청"
$ws.Cells.Item(211,12).Value = "Red"
$ws.Cells.Item(211,13).Value = "반대한다."

# Row 212
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A212:L212").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("M212").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(212,1).Value = 45190.54016662037
$ws.Cells.Item(212,2).Value = "abcdefg76500@naver.com"
$ws.Cells.Item(212,3).Value = "경영대학"
$ws.Cells.Item(212,4).Value = 20232946
$ws.Cells.Item(212,5).Value = "노희망"
$ws.Cells.Item(212,6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(212,7).Value = 0.1
$ws.Cells.Item(212,8).Value = "6:4"
$ws.Cells.Item(212,9).Value = "30분의 1"
$ws.Cells.Item(212,10).Value = "20만호, 69만명"
$ws.Cells.Item(212,11).Value = "전라"
$ws.Cells.Item(212,12).Value = "Red"
$ws.Cells.Item(212,13).Value = "반대한다."

# Row 213
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A213:L213").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("M213").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(213,1).Value = 45190.56663725694
$ws.Cells.Item(213,2).Value = "phyi0915@naver.com"
$ws.Cells.Item(213,3).Value = "경영대학"
$ws.Cells.Item(213,4).Value = 20232967
$ws.Cells.Item(213,5).Value = "백승진"
$ws.Cells.Item(213,6).Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Cells.Item(213,7).Value = 0.5
$ws.Cells.Item(213,8).Value = "3:7"
$ws.Cells.Item(213,9).Value = "15분의 1"
$ws.Cells.Item(213,10).Value = "15만호,  32만명"
$ws.Cells.Item(213,11).Value = "경기"
$ws.Cells.Item(213,12).Value = "Red"
$ws.Cells.Item(213,13).Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# Row 214
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A214:L214").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("M214").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(214,1).Value = 45190.59454880787
$ws.Cells.Item(214,2).Value = "jjedalee@gmail.com"
$ws.Cells.Item(214,3).Value = "심리학과"
$ws.Cells.Item(214,4).Value = 20232134
$ws.Cells.Item(214,5).Value = "정재은"
$ws.Cells.Item(214,6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(214,7).Value = 0.9
$ws.Cells.Item(214,8).Value = "6:4"
$ws.Cells.Item(214,9).Value = "10분의 1"
$ws.Cells.Item(214,10).Value = "20만호, 69만명"
$ws.Cells.Item(214,11).Value = "충청"
$ws.Cells.Item(214,12).Value = "Red"
$ws.Cells.Item(214,13).Value = "반대한다."

# Row 215
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A215:L215").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("M215").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(215,1).Value = 45190.5952000463
$ws.Cells.Item(215,2).Value = "vlxjvos2514@naver.com"
$ws.Cells.Item(215,3).Value = "환경생명공학과"
$ws.Cells.Item(215,4).Value = 20203739
$ws.Cells.Item(215,5).Value = "홍성은"
$ws.Cells.Item(215,6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(215,7).Value = 0.1
$ws.Cells.Item(215,8).Value = "6:4"
$ws.Cells.Item(215,9).Value = "20분의 1"
$ws.Cells.Item(215,10).Value = "20만호, 69만명"
$ws.Cells.Item(215,11).Value = "충청"
$ws.Cells.Item(215,12).Value = "Red"
$ws.Cells.Item(215,13).Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# Row 216
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A216:L216").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("N216").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(216,1).Value = 45190.617646921295
$ws.Cells.Item(216,2).Value = "soc07030@naver.com"
$ws.Cells.Item(216,3).Value = "러시아학과"
$ws.Cells.Item(216,4).Value = 20201723
$ws.Cells.Item(216,5).Value = "윤현수"
$ws.Cells.Item(216,6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(216,7).Value = 0.1
$ws.Cells.Item(216,8).Value = "7:3"
$ws.Cells.Item(216,9).Value = "10분의 1"
$ws.Cells.Item(216,10).Value = "15만호,  32만명"
$ws.Cells.Item(216,11).Value = "평안"
$ws.Cells.Item(216,12).Value = "Black"
$ws.Cells.Item(216,14).Value = "찬성한다."

# Row 217
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A217:L217").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("M217").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(217,1).Value = 45190.62377230324
$ws.Cells.Item(217,2).Value = "hsm3932@naver.com"
$ws.Cells.Item(217,3).Value = "일본학과"
$ws.Cells.Item(217,4).Value = 20221637
$ws.Cells.Item(217,5).Value = "허선민"
$ws.Cells.Item(217,6).Value = "실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다."
$ws.Cells.Item(217,7).Value = 0.7
$ws.Cells.Item(217,8).Value = "4:6"
$ws.Cells.Item(217,9).Value = "15분의 1"
$ws.Cells.Item(217,10).Value = "130만호, 5백만명"
$ws.Cells.Item(217,11).Value = "충청"
$ws.Cells.Item(217,12).Value = "Red"
$ws.Cells.Item(217,13).Value = "모름/무응답"

# Row 218
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A218:L218").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("N218").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(218,1).Value = 45190.663551840276
$ws.Cells.Item(218,2).Value = "jkw1391@naver.com"
$ws.Cells.Item(218,3).Value = "광고홍보학과"
$ws.Cells.Item(218,4).Value = 20212611
$ws.Cells.Item(218,5).Value = "김지원"
$ws.Cells.Item(218,6).Value = "실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다."
$ws.Cells.Item(218,7).Value = 0.1
$ws.Cells.Item(218,8).Value = "6:4"
$ws.Cells.Item(218,9).Value = "15분의 1"
$ws.Cells.Item(218,10).Value = "20만호, 69만명"
$ws.Cells.Item(218,11).Value = "평안"
$ws.Cells.Item(218,12).Value = "Black"
$ws.Cells.Item(218,14).Value = "모름/무응답"

# Row 219
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A219:L219").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("M219").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(219,1).Value = 45190.66730758102
$ws.Cells.Item(219,2).Value = "kimcr0678@gmail.com"
$ws.Cells.Item(219,3).Value = "간호학과"
$ws.Cells.Item(219,4).Value = 20236230
$ws.Cells.Item(219,5).Value = "김채령"
$ws.Cells.Item(219,6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(219,7).Value = 0.1
$ws.Cells.Item(219,8).Value = "5:5"
$ws.Cells.Item(219,9).Value = "10분의 1"
$ws.Cells.Item(219,10).Value = "20만호, 69만명"
$ws.Cells.Item(219,11).Value = "충청"
$ws.Cells.Item(219,12).Value = "Red"
$ws.Cells.Item(219,13).Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# Row 220
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A220:L220").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("N220").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(220,1).Value = 45190.67234148148
$ws.Cells.Item(220,2).Value = "wnyuna04@gmail.com"
$ws.Cells.Item(220,3).Value = "간호학과"
$ws.Cells.Item(220,4).Value = 20236294
$ws.Cells.Item(220,5).Value = "주윤아"
$ws.Cells.Item(220,6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(220,7).Value = 0.1
$ws.Cells.Item(220,8).Value = "6:4"
$ws.Cells.Item(220,9).Value = "20분의 1"
$ws.Cells.Item(220,10).Value = "20만호, 69만명"
$ws.Cells.Item(220,11).Value = "충청"
$ws.Cells.Item(220,12).Value = "Black"
$ws.Cells.Item(220,14).Value = "찬성한다."

# Row 221
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A221:L221").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("M221").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(221,1).Value = 45190.71674582176
$ws.Cells.Item(221,2).Value = "tvkim0122@gmail.com"
$ws.Cells.Item(221,3).Value = "콘텐츠IT"
$ws.Cells.Item(221,4).Value = 20203314
$ws.Cells.Item(221,5).Value = "김태빈"
$ws.Cells.Item(221,6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(221,7).Value = 0.1
$ws.Cells.Item(221,8).Value = "6:4"
$ws.Cells.Item(221,9).Value = "20분의 1"
$ws.Cells.Item(221,10).Value = "20만호, 69만명"
$ws.Cells.Item(221,11).Value = "충청"
$ws.Cells.Item(221,12).Value = "Red"
$ws.Cells.Item(221,13).Value = "모름/무응답"

# Row 222
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A222:L222").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("M222").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(222,1).Value = 45190.72814460648
$ws.Cells.Item(222,2).Value = "pgw0814@gmail.com"
$ws.Cells.Item(222,3).Value = "화학과"
$ws.Cells.Item(222,4).Value = 20203409
$ws.Cells.Item(222,5).Value = "박지우"
$ws.Cells.Item(222,6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(222,7).Value = 0.1
$ws.Cells.Item(222,8).Value = "6:4"
$ws.Cells.Item(222,9).Value = "20분의 1"
$ws.Cells.Item(222,10).Value = "20만호, 69만명"
$ws.Cells.Item(222,11).Value = "충청"
$ws.Cells.Item(222,12).Value = "Red"
$ws.Cells.Item(222,13).Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# Row 223
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A223:L223").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("N223").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(223,1).Value = 45190.755936331014
$ws.Cells.Item(223,2).Value = "aoa0226@naver.com"
$ws.Cells.Item(223,3).Value = "일본학과"
$ws.Cells.Item(223,4).Value = 20201601
$ws.Cells.Item(223,5).Value = "강문원"
$ws.Cells.Item(223,6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(223,7).Value = 0.1
$ws.Cells.Item(223,8).Value = "6:4"
$ws.Cells.Item(223,9).Value = "20분의 1"
$ws.Cells.Item(223,10).Value = "20만호, 69만명"
$ws.Cells.Item(223,11).Value = "충청"
$ws.Cells.Item(223,12).Value = "Black"
$ws.Cells.Item(223,14).Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# Row 224
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A224:L224").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("M224").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(224,1).Value = 45190.755951956016
$ws.Cells.Item(224,2).Value = "0218mun@naver.com"
$ws.Cells.Item(224,3).Value = "바이오매디컬"
$ws.Cells.Item(224,4).Value = 20193601
$ws.Cells.Item(224,5).Value = "강문희"
$ws.Cells.Item(224,6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(224,7).Value = 0.1
$ws.Cells.Item(224,8).Value = "6:4"
$ws.Cells.Item(224,9).Value = "20분의 1"
$ws.Cells.Item(224,10).Value = "20만호, 69만명"
$ws.Cells.Item(224,11).Value = "충청"
$ws.Cells.Item(224,12).Value = "Red"
$ws.Cells.Item(224,13).Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# Row 225
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A225:L225").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("N225").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(225,1).Value = 45190.78861594907
$ws.Cells.Item(225,2).Value = "vivace3990@gmail.com"
$ws.Cells.Item(225,3).Value = "일본학과"
$ws.Cells.Item(225,4).Value = 20171605
$ws.Cells.Item(225,5).Value = "김정현"
$ws.Cells.Item(225,6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(225,7).Value = 0.9
$ws.Cells.Item(225,8).Value = "6:4"
$ws.Cells.Item(225,9).Value = "20분의 1"
$ws.Cells.Item(225,10).Value = "20만호, 69만명"
$ws.Cells.Item(225,11).Value = "충청"
$ws.Cells.Item(225,12).Value = "Black"
$ws.Cells.Item(225,14).Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# Row 226
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A226:L226").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("N226").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(226,1).Value = 45190.79722414352
$ws.Cells.Item(226,2).Value = "ily0302@naver.com"
$ws.Cells.Item(226,3).Value = "러시아"
$ws.Cells.Item(226,4).Value = 20221730
$ws.Cells.Item(226,5).Value = "최요원"
$ws.Cells.Item(226,6).Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Cells.Item(226,7).Value = 0.5
$ws.Cells.Item(226,8).Value = "4:6"
$ws.Cells.Item(226,9).Value = "20분의 1"
$ws.Cells.Item(226,10).Value = "20만호, 69만명"
$ws.Cells.Item(226,11).Value = "전라"
$ws.Cells.Item(226,12).Value = "Black"
$ws.Cells.Item(226,14).Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# Row 227
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A227:L227").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("M227").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(227,1).Value = 45190.82719306713
$ws.Cells.Item(227,2).Value = "jieunpark0408@naver.com"
$ws.Cells.Item(227,3).Value = "인문학부"
$ws.Cells.Item(227,4).Value = 20231043
$ws.Cells.Item(227,5).Value = "박지은"
$ws.Cells.Item(227,6).Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Cells.Item(227,7).Value = 0.3
$ws.Cells.Item(227,8).Value = "7:3"
$ws.Cells.Item(227,9).Value = "15분의 1"
$ws.Cells.Item(227,10).Value = "20만호, 69만명"
$ws.Cells.Item(227,11).Value = "전라"
$ws.Cells.Item(227,12).Value = "Red"
$ws.Cells.Item(227,13).Value = "반대한다."

# Row 228
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A228:L228").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("N228").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(228,1).Value = 45190.839383391205
$ws.Cells.Item(228,2).Value = "gksgh2311@naver.com"
$ws.Cells.Item(228,3).Value = "경제학과"
$ws.Cells.Item(228,4).Value = 20232842
$ws.Cells.Item(228,5).Value = "전한호"
$ws.Cells.Item(228,6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(228,7).Value = 0.5
$ws.Cells.Item(228,8).Value = "6:4"
$ws.Cells.Item(228,9).Value = "10분의 1"
$ws.Cells.Item(228,10).Value = "20만호, 69만명"
$ws.Cells.Item(228,11).Value = "경기"
$ws.Cells.Item(228,12).Value = "Black"
$ws.Cells.Item(228,14).Value = "모름/무응답"

# Row 229
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A229:L229").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("N229").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(229,1).Value = 45190.85034533565
$ws.Cells.Item(229,2).Value = "parksiwoo1214@naver.com"
$ws.Cells.Item(229,3).Value = "데이터사이언스학부"
$ws.Cells.Item(229,4).Value = 20193219
$ws.Cells.Item(229,5).Value = "박시우"
$ws.Cells.Item(229,6).Value = "‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."
$ws.Cells.Item(229,7).Value = 0.1
$ws.Cells.Item(229,8).Value = "7:3"
$ws.Cells.Item(229,9).Value = "15분의 1"
$ws.Cells.Item(229,10).Value = "44만호, 153만명"
$ws.Cells.Item(229,11).Value = "경기"
$ws.Cells.Item(229,12).Value = "Black"
$ws.Cells.Item(229,14).Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# Row 230
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A230:L230").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("M230").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(230,1).Value = 45190.85991384259
$ws.Cells.Item(230,2).Value = "kns2266@naver.com"
$ws.Cells.Item(230,3).Value = "사학과"
$ws.Cells.Item(230,4).Value = 20201012
$ws.Cells.Item(230,5).Value = "김남석"
$ws.Cells.Item(230,6).Value = "‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."
$ws.Cells.Item(230,7).Value = 0.7
$ws.Cells.Item(230,8).Value = "4:6"
$ws.Cells.Item(230,9).Value = "15분의 1"
$ws.Cells.Item(230,10).Value = "130만호, 5백만명"
$ws.Cells.Item(230,11).Value = "평안"
$ws.Cells.Item(230,12).Value = "Red"
$ws.Cells.Item(230,13).Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# Row 231
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A231:L231").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("N231").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(231,1).Value = 45190.866045717594
$ws.Cells.Item(231,2).Value = "m.kwak2018@gmail.com"
$ws.Cells.Item(231,3).Value = "사회학과 "
$ws.Cells.Item(231,4).Value = 20182202
$ws.Cells.Item(231,5).Value = "곽민수"
$ws.Cells.Item(231,6).Value = "‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."
$ws.Cells.Item(231,7).Value = 0.1
$ws.Cells.Item(231,8).Value = "5:5"
$ws.Cells.Item(231,9).Value = "10분의 1"
$ws.Cells.Item(231,10).Value = "130만호, 5백만명"
$ws.Cells.Item(231,11).Value = "전라"
$ws.Cells.Item(231,12).Value = "Black"
$ws.Cells.Item(231,14).Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# Row 232
$ws.Range("A203:L203").Copy() | Out-Null
$ws.Range("A232:L232").PasteSpecial(-4122) | Out-Null
$ws.Range("N203").Copy() | Out-Null
$ws.Range("M232").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(232,1).Value = 45190.89194008102
$ws.Cells.Item(232,2).Value = "popoq2004@naver.com"
$ws.Cells.Item(232,3).Value = "간호학과"
$ws.Cells.Item(232,4).Value = 20236239
$ws.Cells.Item(232,5).Value = "박신비"
$ws.Cells.Item(232,6).Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Cells.Item(232,7).Value = 0.3
$ws.Cells.Item(232,8).Value = "4:6"
$ws.Cells.Item(232,9).Value = "15분의 1"
$ws.Cells.Item(232,10).Value = "44만호, 153만명"
$ws.Cells.Item(232,11).Value = "전라"
$ws.Cells.Item(232,12).Value = "Red"
$ws.Cells.Item(232,13).Value = "반대한다."

$excel.CutCopyMode = 0
